$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 5 held "5. Explore AI for All" (text in A5/B5 + hyperlinks on B5
# and C5). That whole row's content is being removed, while rows 6 and
# 7 (cybersecurity / AI security) keep their row numbers and formatting.
# ---------------------------------------------------------------------
$ws.Range("A5:C5").ClearContents()

# Deleting via a Range wipes the sheet's whole Hyperlinks collection in
# this host, so capture + restore every remaining link explicitly.
$ws.Hyperlinks.Delete()

function Add-Link($addr, $target, $text) {
    $ws.Hyperlinks.Add($ws.Range($addr), $target, "", "", $target) | Out-Null
    $ws.Range($addr).Value = $text
}

Add-Link "B1" "https://learn.microsoft.com/en-us/training/modules/explore-ai-basics/" "1. Explore AI basics"
Add-Link "B2" "https://learn.microsoft.com/en-us/training/modules/explore-generative-ai/" "2. Explore Generative AI"
Add-Link "B3" "https://learn.microsoft.com/en-us/training/modules/explore-internet-search-beyond/" "3. Explore internet search and beyond"
Add-Link "B4" "https://learn.microsoft.com/en-us/training/modules/responsible-ai/" "4. Explore responsible AI"
Add-Link "B6" "https://learn.microsoft.com/en-us/training/paths/describe-basic-concepts-of-cybersecurity/" "6. Describe the concepts of cybersecurity "
Add-Link "B7" "https://learn.microsoft.com/en-us/training/paths/ai-security-fundamentals/" "7. AI security fundamentals"

$ws.Hyperlinks.Add($ws.Range("C1"), "https://share.articulate.com/DYPFmXtdrG2phDVMk34SF") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://share.articulate.com/9cGYU-Jl9BFbmyWtP0cN9") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://share.articulate.com/ZM2j5uOMQAttXsvQAICdn") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://share.articulate.com/Bh76w2pleoYOz-Te5B4H9") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://share.articulate.com/xmBvMFOeW7vUMXhLP6Ri7") | Out-Null

Add-Link "C7" "https://share.articulate.com/fc9LKnz5xeESiGpJ4Tydt" "https://share.articulate.com/h0wMTOonB6tfpw1TJQ7kM`nhttps://share.articulate.com/zKp2T-P9fUMKimzPN3hUG`nhttps://share.articulate.com/fc9LKnz5xeESiGpJ4Tydt"

# ---------------------------------------------------------------------
# View / formatting tweaks that came along with the re-save.
# ---------------------------------------------------------------------
$ws.Range("A1:C7").RowHeight = 16.5
$ws.Range("A7:C7").RowHeight = 17.25

$ws.Range("A5:G5").Select()
